$wb = $excel.ActiveWorkbook

# Per-sheet new coefficient values for rows 2-12 (after deleting the bike_lane_share_res row)
$sheetData = @{}
$sheetData[1] = @(1.681349514153613,-13.63196077028263,-24.89251717533051,0.834730103705339,-0.0004569554875543484,52.05863061667264,31.4366590255565,-1.522656592477018,-1.109443211756551,-0.2756264011924068,-73.09573345049176)
$sheetData[2] = @(0.658882427413374,-2.773615115779819,-0.8435543175151281,-0.2738163954625193,-0.0004419014771841481,42.23026168238317,40.62258080672893,-0.1803292559782821,0.06042693051862658,0.06747092731359661,6.457818355140148)
$sheetData[3] = @(64.85441897555307,673.3095553097511,463.9712135788503,-0.9821200187678065,0.0009351434254842284,-159.7843447472493,-143.1234889330963,9.456611981140171,30.42545868015272,-8.783859022917333,681.961570128486)
$sheetData[4] = @(2.998129462739111,-31.24473533092061,-24.4023501245104,0.5719126912802819,-0.001974145998904666,24.52393982145726,-4.609423762209033,-0.8184787431172765,-0.4142520150779861,0.4470985364340567,166.3392089275536)
$sheetData[5] = @(3.396393335818191,-3.858322073310496,-55.49471631395748,1.776882263254514,-0.005565623695234459,-18.96062524137571,52.69738733139506,-0.3431837337131698,-1.70601222687703,0.303967575327228,74.32876685645327)
$sheetData[6] = @(-5.325575769836185,94.56710763192936,450.4820752125969,-0.5306133458539588,0.001369270690653088,4.157807000605146,14.10257701937832,1.77897001347797,15.80932995840931,-2.665729464034732,-17.87921330643792)
$sheetData[7] = @(0.8789658349713001,2.79655682585657,-13.09194798707108,1.687175310859491,-0.003945878204174599,-61.95195079936207,50.38867387890552,-0.1252204978490603,-0.7110796044299867,0.04116524525667423,25.35929025959962)
$sheetData[8] = @(4.277389175381323,-45.1030158575832,-40.15170553404363,-0.6564526932867625,-0.0007477211087680367,40.65656492730125,-9.755779280099972,-1.550976961570021,-0.9294865186935963,0.4670717883833742,188.6811255083308)
$sheetData[9] = @(4.438285736636729,-8.174538076944202,-75.22703510367046,0.04196943477477388,-0.001102580471154237,-4.448127237952242,29.47042645170501,-0.5295007773738343,-2.178962074203412,0.3626933692778324,89.2101102422273)
$sheetData[10] = @(-3.848383160623807,82.43253829002171,251.8958517568991,2.466000822373,-0.006153005788394014,64.96550582232726,72.5726749617083,2.245318913546326,9.096865863234481,-1.425994517181972,-92.71203958610877)
$sheetData[11] = @(2.445403585617377,-13.98862900215779,-45.8100077309733,-0.6432388293135234,0.0002509407630924801,-22.68468170129671,80.20000948818009,-1.819494061037485,-2.227789346556183,-0.236319818397326,-57.73399213483167)
$sheetData[12] = @(-5.325575769836185,94.56710763192936,450.4820752125969,-0.5306133458539588,0.001369270690653088,4.157807000605146,14.10257701937832,1.77897001347797,15.80932995840931,-2.665729464034732,-17.87921330643792)
$sheetData[13] = @(1.213396391457235,2.854949600910603,1.640895628936185,0.8145724096418774,-0.003028113414281689,53.59796974142732,22.91157256676626,0.5210348359510353,0.6283168540139845,0.2813121628130568,72.52619575950757)
$sheetData[14] = @(5.082615964617212,-0.5615654751998953,-87.64945419414948,0.05651028117492018,-0.001014799657284568,-2.389415753843394,32.3506000080628,-0.4047759855066269,-2.495170878277076,0.3284886399237941,73.11485463996075)
$sheetData[15] = @(109.8837683654954,-860.3620878216009,1471.331969527296,0.1486737854178632,0.001564936471318038,239.0001488810807,78.99682190879668,-102.7385633038138,35.18701064087962,-32.73573380482163,-413.4004118632636)
$sheetData[16] = @(5.082615964617212,-0.5615654751998953,-87.64945419414948,0.05651028117492018,-0.001014799657284568,-2.389415753843394,32.3506000080628,-0.4047759855066269,-2.495170878277076,0.3284886399237941,73.11485463996075)
$sheetData[17] = @(0.7342457431960945,8.511683820079661,9.622981079129985,2.003004533750894,-0.006272507897412123,-22.92226610390486,49.9785869365954,0.7262239587855617,0.5728190467685725,0.269631431170096,90.8967625199117)
$sheetData[18] = @(20.37756816783126,-236.5299341595297,-314.1550276709117,0.5872185148438704,-0.001811060301938866,75.4394643237372,16.10051867252142,-15.12407246289665,-12.83995605448061,-0.7093011909765259,77.92728140997099)
$sheetData[19] = @(0.2628345745311442,-15.17943046457873,-0.4927659539268348,-0.3269791275883671,-0.0005157582627849707,-32.94152429233745,45.85934378505956,-0.4803939397121404,-0.3267394996061705,0.1798507109815961,59.82164893913011)
$sheetData[20] = @(-4.362506641510457,88.37166187270347,269.7325549594341,2.40195532482142,-0.005889500845754009,53.16824153110007,63.77296577777722,2.788905257902875,9.84921812782987,-1.369707504611111,-61.26565827397297)
$sheetData[21] = @(-0.2766408849400861,-9.596981121963395,22.82413883257351,0.591360752164638,0.002025502279504968,-134.8055612727211,44.86369995308406,-2.306010101399611,-0.7426329335980055,-0.943759408847179,-121.630106864127)
$sheetData[22] = @(3.00745182993669,-36.78832077086213,-62.76870180196771,-1.057497401348696,0.0007107069179198471,-67.99580763932359,92.76353532714415,-3.109766981297661,-3.374228586202158,-0.2950312056941144,-49.92343089263801)
$sheetData[23] = @(0.7342457431960945,8.511683820079661,9.622981079129985,2.003004533750894,-0.006272507897412123,-22.92226610390486,49.9785869365954,0.7262239587855617,0.5728190467685725,0.269631431170096,90.8967625199117)
$sheetData[24] = @(28.91227899216383,241.8145701889192,360.1621361724408,0.1608472210046853,-0.004604112461823854,-87.77500597807806,-89.97821669512641,2.893927280758109,18.95029987897884,-3.816946878147865,662.256758245063)
$sheetData[25] = @(2.09686796752406,-3.995463131682335,-37.4005485404906,-0.1455430610778663,-0.0009836119598913545,-25.93417105720937,76.52235349175157,-1.067637252603317,-1.704936092651209,-0.1200435764807106,-32.98296813468744)
$sheetData[26] = @(2.803157531852992,4.665553439860013,288.2705624587002,-2.67886162413695,0.0009024169688282371,-44.81060387315964,-46.66502791583949,-0.4619825929502234,10.92513034215301,-1.148765866404534,394.7899138688308)
$sheetData[27] = @(28.91227899216383,241.8145701889192,360.1621361724408,0.1608472210046853,-0.004604112461823854,-87.77500597807806,-89.97821669512641,2.893927280758109,18.95029987897884,-3.816946878147865,662.256758245063)
$sheetData[28] = @(0.4506132996183717,3.6792716843753,4.815573864331313,0.235381208224938,-0.0004090224003885016,44.24841157015709,15.87203290271702,0.2983252889773942,0.4793245283644911,0.0999695373936435,18.6856947578444)
$sheetData[29] = @(3.272318707646808,-3.979999676640205,-53.83925303658799,2.351497380874829,-0.00703426475710815,-22.25919584180419,61.51609650875901,-0.427206273854822,-1.740452013809509,0.2700959430260926,64.76103252724472)
$sheetData[30] = @(10.94285616068795,-131.7751585596457,-208.8661858221026,-0.1498265269254935,-0.0006867771747092568,10.27020635311345,50.83638191064065,-8.286441527054638,-8.631691263166397,-0.1677216937376922,20.39204726849056)

# New sheet names in rId1..rId30 / sheet1..sheet30 order
$newNames = @("summ51438189","summ51570461","summ51702614","summ51843257","summ51977212","summ52112911","summ52254624","summ52396434","summ52546195","summ52683914","summ52862434","summ53005262","summ53145519","summ53288485","summ53426198","summ53568586","summ53709208","summ53849972","summ53987626","summ54163449","summ54322762","summ54474293","summ54617185","summ54771765","summ54910781","summ55148621","summ55286191","summ55422968","summ55565166","summ55705493")

for ($i = 1; $i -le 30; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Delete row 11 (bike_lane_share_res), shifting Commute_Trip/Age up to rows 11/12
    $ws.Rows(11).Delete()

    # Write the new coefficient values into B2:B12
    $vals = $sheetData[$i]
    for ($r = 0; $r -lt 11; $r++) {
        $ws.Cells.Item($r + 2, 2).Value2 = $vals[$r]
    }

    # Rename the sheet
    $ws.Name = $newNames[$i - 1]
}

Write-Output "Done"